# "Generate Report for Handback": refresh the handback timestamps for the
# cf1bd38b-3131-4f9f-a372-8fd37e9375a9 row across the Overview / zh-cn / de-de
# sheets to reflect a newly (re)generated handback report.

$wb = $excel.ActiveWorkbook

$wsOverview = $wb.Worksheets.Item("Overview")
$wsZhCn     = $wb.Worksheets.Item("zh-cn")
$wsDeDe     = $wb.Worksheets.Item("de-de")

# Overview!G2 - "Latest HO Xliff Generate Date" for cf1bd38b-...
$wsOverview.Range("G2").Value = "2016-10-17 15:36:13"

# zh-cn!H2 - "Correspond Handoff Datetime" for cf1bd38b-...
$wsZhCn.Range("H2").Value = "2016-10-17 15:35:50"

# zh-cn!K2 - "Correspond Handback DateTime" for cf1bd38b-...
$wsZhCn.Range("K2").Value = "2016-10-17 15:36:59"

# de-de!H2 - "Correspond Handoff Datetime" for cf1bd38b-...
$wsDeDe.Range("H2").Value = "2016-10-17 15:36:13"

# de-de!K2 - "Correspond Handback DateTime" for cf1bd38b-...
$wsDeDe.Range("K2").Value = "2016-10-17 15:37:37"
